$d = $word.ActiveDocument

# --- 1. Replace the placeholder heading with Myson's real second-pass writeup ---
$rng = $d.Content
$found1 = $rng.Find.Execute("******MYSON’S STUFF GO HERE*****", $true, $false, $false, $false, $false, $true, 1, $false, "In the second pass, Myson took the intermediate file from John’s first pass in order to generate the object code and finally the object program. Our approach to the design was difficult at first because designing top to bottom while also designing bottom up is hard when there are a lot of dependencies. Initially in the secondpass object, Myson stepped through each line in the intermediate file and parsed that file that now included the LOC for each line which is pivotal in creating the object code. This firstpassoutput object used the split() functionality that John discovered to be useful in the first pass. A lot of the second pass depends on StringBuilders for building upon the object program and the parsed intermediate text file to generate the object code. While stepping through each line of the intermediate text file, the second pass will check to see what format the instruction is and generate the object code accordingly. The object directly deals with format 1 and 2 because there is not a lot of calculations to do there. If the format is 3 or 4, then that line will be passed to the objcodegen object that will deal with all the details of PC/Base relative, target addresses and disp calculations. This object also deals with immediate, indexed and indirect addressing modes. There is a lot of code here that definitely needs more thorough debugging to ensure it works properly. Some of the difficulty with the object code generator object came about when actually computing the hexadecimal values. Myson initially tried to do some calculations that were not accurate but once discovering all of the built in Java functions for bytes and hexadecimal formatting everything worked fine. After generating the object code, the second pass also handles instructions that do not generate object code and other special instruction such as BYTE and WORD. The pass finally ends with appending the object code to the text record in the correct formats. One tricky thing that John was able to catch were the half byte representations in the records.", 2)
Write-Host "Step1 replace MYSON placeholder: $found1"

# --- 3. Move the _GoBack bookmark to sit right after the hex-formatting sentence ---
$hadBookmark = $d.Bookmarks.Exists("_GoBack")
Write-Host "Had _GoBack before: $hadBookmark"
if ($hadBookmark) {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
}
$rng3 = $d.Content
$found2 = $rng3.Find.Execute("hexadecimal formatting everything worked fine. ", $true)
Write-Host "Step3 find bookmark anchor: $found2"
if ($found2) {
    $bmRange = $d.Range($rng3.End, $rng3.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# --- 4. Insert the new sentence about incomplete/inaccurate object programs ---
$rng4 = $d.Content
$found3 = $rng4.Find.Execute("contained literals as well. We utilized github to post", $true, $false, $false, $false, $false, $true, 1, $false, "contained literals as well. Right now, each example provided will provide an object program but it will not be accurate because we were not able to complete the project. We utilized github to post", 2)
Write-Host "Step4 insert new sentence: $found3"
